$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.021625995635986
$ws.Range("B1").Value = 2.758477449417114
$ws.Range("C1").Value = 6.743141651153564
$ws.Range("D1").Value = 4.099586963653564
$ws.Range("E1").Value = 1.432390093803406
